$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Honduras" / "Bulgaria" country names so Bulgaria now precedes Honduras
$ws.Cells.Item(79, 1).Value = "Bulgaria"
$ws.Cells.Item(80, 1).Value = "Honduras"

# Row 68 (Hungria) - updated daily stats
$ws.Cells.Item(68, 2).Value = 3284
$ws.Cells.Item(68, 3).Value = 21
$ws.Cells.Item(68, 4).Value = 958
$ws.Cells.Item(68, 5).Value = 1905
$ws.Cells.Item(68, 6).Value = 42
$ws.Cells.Item(68, 7).Value = 8
$ws.Cells.Item(68, 8).Value = 421

# Row 75 (Uzbekistan) - updated daily stats
$ws.Cells.Item(75, 2).Value = 2453
$ws.Cells.Item(75, 3).Value = 35
$ws.Cells.Item(75, 5).Value = 562

# Row 79 (now Bulgaria) - updated daily stats
$ws.Cells.Item(79, 2).Value = 1981
$ws.Cells.Item(79, 3).Value = 16
$ws.Cells.Item(79, 4).Value = 461
$ws.Cells.Item(79, 5).Value = 1429
$ws.Cells.Item(79, 6).Value = 58
$ws.Cells.Item(79, 8).Value = 91

# Row 80 (now Honduras) - updated daily stats
$ws.Cells.Item(80, 2).Value = 1972
$ws.Cells.Item(80, 3).Value = 142
$ws.Cells.Item(80, 4).Value = 203
$ws.Cells.Item(80, 5).Value = 1661
$ws.Cells.Item(80, 6).Value = 10
$ws.Cells.Item(80, 8).Value = 108
